# Update gh-pages to output generated at 456a3b4
# Apply numeric updates to column F ("浏览数"/view-count style counters)
# across the four worksheets: 展览, 演出, 本地生活, 全部类型.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet
$wsExhibit.Range("F12").Value = 73
$wsExhibit.Range("F15").Value = 24
$wsExhibit.Range("F16").Value = 685
$wsExhibit.Range("F17").Value = 176
$wsExhibit.Range("F20").Value = 8275
$wsExhibit.Range("F47").Value = 198

# 演出 (Performance) sheet
$wsShow.Range("F7").Value = 18
$wsShow.Range("F18").Value = 306

# 本地生活 (Local Life) sheet
$wsLocal.Range("F5").Value = 147

# 全部类型 (All Types) sheet
$wsAll.Range("F6").Value = 147
$wsAll.Range("F14").Value = 73
$wsAll.Range("F15").Value = 24
$wsAll.Range("F16").Value = 685
$wsAll.Range("F18").Value = 176
$wsAll.Range("F22").Value = 8275
$wsAll.Range("F47").Value = 198
$wsAll.Range("F49").Value = 306
